# Ran code for averaged intensities on spiral schemes
# Inserts 3 new rotation schemes ("Spiral-90deg-...") into the averaged
# intensity table, pushing the existing "NoRotation/Rotation/HexGrid"
# schemes down and appending the table's final scheme to new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the existing scheme labels in column B (rows 10-16) down to
#     make room for the 3 new spiral schemes that now land at rows 11-13 ---
$ws.Cells.Item(10, 2).Value = "Gaussian-Quadrature"
$ws.Cells.Item(11, 2).Value = "Spiral-90deg-10rot-5space"
$ws.Cells.Item(12, 2).Value = "Spiral-90deg-15rot-5space"
$ws.Cells.Item(13, 2).Value = "Spiral-90deg-10rot-3space"
$ws.Cells.Item(14, 2).Value = "NoRotation-tilt60deg"
$ws.Cells.Item(15, 2).Value = "Rotation-NoTilt"
$ws.Cells.Item(16, 2).Value = "Rotation-60detTilt"

# --- Append the 3 rows that fell off the end of the table (rows 17-19) ---
$newRows = @(
    @{ Row = 17; Idx = 15; Name = "HexGrid-90degTilt5degRes" },
    @{ Row = 18; Idx = 16; Name = "HexGrid-90degTilt22p5degRes" },
    @{ Row = 19; Idx = 17; Name = "HexGrid-60degTilt5degRes" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Idx
    $ws.Cells.Item($row, 2).Value = $r.Name

    # Averaged intensities (columns C:P) for the row are all 1
    for ($col = 3; $col -le 16; $col++) {
        $ws.Cells.Item($row, $col).Value = 1
    }

    # Column A uses the same bold/centered/bordered style as the rest of
    # the index column - copy it from the row above instead of rebuilding
    # it property-by-property so it resolves to the existing style.
    $ws.Cells.Item($row - 1, 1).Copy() | Out-Null
    $ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null
}
